$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6908.6816
$ws.Range("J19").Value = 8684.117
$ws.Range("L19").Value = 8684.117
$ws.Range("N19").Value = -9034.117

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1471.1904
$ws.Range("I98").Value = 1471.1904
$ws.Range("K98").Value = 1471.1904
$ws.Range("M98").Value = 26.80960000000005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4959.5454
$ws.Range("I113").Value = 5341.3335
$ws.Range("J113").Value = 4141.4287
$ws.Range("K113").Value = 5341.3335
$ws.Range("L113").Value = 4141.4287
$ws.Range("M113").Value = -2087.3335
$ws.Range("N113").Value = -10649.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1471.1904
$ws.Range("I122").Value = 1471.1904
$ws.Range("K122").Value = 4413.5712
$ws.Range("M122").Value = -1963.5712

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1011949.94
$ws.Range("I129").Value = 2033.8
$ws.Range("J129").Value = 11111111
$ws.Range("K129").Value = 6101.4
$ws.Range("L129").Value = 33333333
$ws.Range("M129").Value = -1101.4
$ws.Range("N129").Value = -33343333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5023.0625
$ws.Range("I132").Value = 3977.851
$ws.Range("J132").Value = 7912.7646
$ws.Range("K132").Value = 11933.553
$ws.Range("L132").Value = 23738.2938
$ws.Range("M132").Value = -9403.553
$ws.Range("N132").Value = -28798.2938

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5152.8667
$ws.Range("I137").Value = 5817.6665
$ws.Range("J137").Value = 4709.6665
$ws.Range("K137").Value = 17452.9995
$ws.Range("L137").Value = 14128.9995
$ws.Range("M137").Value = -14902.9995
$ws.Range("N137").Value = -19228.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 82058380
$ws.Range("I138").Value = 111116600
$ws.Range("J138").Value = 16677416
$ws.Range("K138").Value = 333349800
$ws.Range("L138").Value = 50032248
$ws.Range("M138").Value = -333344660
$ws.Range("N138").Value = -50042528

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3303.6
$ws.Range("I61").Value = 2136.2856
$ws.Range("J61").Value = 6027.3335
$ws.Range("K61").Value = 2136.2856
$ws.Range("L61").Value = 6027.3335
$ws.Range("M61").Value = -1924.2856
$ws.Range("N61").Value = -6451.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3303.6
$ws.Range("I136").Value = 2136.2856
$ws.Range("J136").Value = 6027.3335
$ws.Range("K136").Value = 6408.8568
$ws.Range("L136").Value = 18082.0005
$ws.Range("M136").Value = -3858.8568
$ws.Range("N136").Value = -23182.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2391.818
$ws.Range("I86").Value = 1833
$ws.Range("J86").Value = 3882
$ws.Range("K86").Value = 1833
$ws.Range("L86").Value = 3882
$ws.Range("M86").Value = -710
$ws.Range("N86").Value = -6128

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2391.818
$ws.Range("I89").Value = 1833
$ws.Range("J89").Value = 3882
$ws.Range("K89").Value = 9165
$ws.Range("L89").Value = 19410
$ws.Range("M89").Value = -3549
$ws.Range("N89").Value = -30642

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2840.5117
$ws.Range("I31").Value = 2623.2144
$ws.Range("J31").Value = 2945.4138
$ws.Range("K31").Value = 2623.2144
$ws.Range("L31").Value = 2945.4138
$ws.Range("M31").Value = -2328.2144
$ws.Range("N31").Value = -3535.4138

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2840.5117
$ws.Range("I34").Value = 2623.2144
$ws.Range("J34").Value = 2945.4138
$ws.Range("K34").Value = 2623.2144
$ws.Range("L34").Value = 2945.4138
$ws.Range("M34").Value = -2421.2144
$ws.Range("N34").Value = -3349.4138

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9132.706
$ws.Range("I99").Value = 6288.7
$ws.Range("J99").Value = 13195.571
$ws.Range("K99").Value = 6288.7
$ws.Range("L99").Value = 13195.571
$ws.Range("M99").Value = -4790.7
$ws.Range("N99").Value = -16191.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1520.425
$ws.Range("I122").Value = 1384.2
$ws.Range("J122").Value = 1929.1
$ws.Range("K122").Value = 4152.6
$ws.Range("L122").Value = 5787.299999999999
$ws.Range("M122").Value = -1702.6
$ws.Range("N122").Value = -10687.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9132.706
$ws.Range("I126").Value = 6288.7
$ws.Range("J126").Value = 13195.571
$ws.Range("K126").Value = 18866.1
$ws.Range("L126").Value = 39586.713
$ws.Range("M126").Value = -16396.1
$ws.Range("N126").Value = -44526.713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 174750
$ws.Range("I141").Value = 50000
$ws.Range("J141").Value = 216333.33
$ws.Range("K141").Value = 50000
$ws.Range("L141").Value = 216333.33
$ws.Range("M141").Value = -44820
$ws.Range("N141").Value = -226693.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 580.53845
$ws.Range("I5").Value = 594.3
$ws.Range("J5").Value = 534.6667
$ws.Range("K5").Value = 1782.9
$ws.Range("L5").Value = 1604.0001
$ws.Range("M5").Value = -1670.9
$ws.Range("N5").Value = -1828.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 25071620
$ws.Range("J32").Value = 25071620
$ws.Range("L32").Value = 75214860
$ws.Range("N32").Value = -75215426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 12030.3
$ws.Range("I68").Value = 2100
$ws.Range("J68").Value = 51751.5
$ws.Range("K68").Value = 6300
$ws.Range("L68").Value = 155254.5
$ws.Range("M68").Value = -5489
$ws.Range("N68").Value = -156876.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 12030.3
$ws.Range("I71").Value = 2100
$ws.Range("J71").Value = 51751.5
$ws.Range("K71").Value = 18900
$ws.Range("L71").Value = 465763.5
$ws.Range("M71").Value = -14844
$ws.Range("N71").Value = -473875.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 6669176.5
$ws.Range("I121").Value = 55556052
$ws.Range("K121").Value = 166668156
$ws.Range("M121").Value = -166666846

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8271083
$ws.Range("J131").Value = 5563014
$ws.Range("L131").Value = 16689042
$ws.Range("N131").Value = -16699122

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 580.53845
$ws.Range("I135").Value = 594.3
$ws.Range("J135").Value = 534.6667
$ws.Range("K135").Value = 5348.7
$ws.Range("L135").Value = 4812.0003
$ws.Range("M135").Value = -2813.7
$ws.Range("N135").Value = -9882.0003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6263.222
$ws.Range("I137").Value = 4770.4
$ws.Range("K137").Value = 14311.2
$ws.Range("M137").Value = -9211.199999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 35269.332
$ws.Range("I141").Value = 35269.332
$ws.Range("K141").Value = 105807.996
$ws.Range("M141").Value = -100627.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 122000
$ws.Range("J42").Value = 122000
$ws.Range("L42").Value = 122000
$ws.Range("N42").Value = -122970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 122000
$ws.Range("J115").Value = 122000
$ws.Range("L115").Value = 122000
$ws.Range("N115").Value = -124350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2772.2144
$ws.Range("I122").Value = 3456
$ws.Range("J122").Value = 2498.7
$ws.Range("K122").Value = 10368
$ws.Range("L122").Value = 7496.099999999999
$ws.Range("M122").Value = -7918
$ws.Range("N122").Value = -12396.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4584.7827
$ws.Range("I132").Value = 4247.8335
$ws.Range("K132").Value = 12743.5005
$ws.Range("M132").Value = -10213.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10950.583
$ws.Range("I61").Value = 13489
$ws.Range("K61").Value = 13489
$ws.Range("M61").Value = -13287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 10950.583
$ws.Range("I113").Value = 13489
$ws.Range("K113").Value = 13489
$ws.Range("M113").Value = -11319

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 80235
$ws.Range("J54").Value = 80235
$ws.Range("L54").Value = 80235
$ws.Range("N54").Value = -81275

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4813.7456
$ws.Range("I122").Value = 3174.98
$ws.Range("J122").Value = 13918
$ws.Range("K122").Value = 9524.940000000001
$ws.Range("L122").Value = 41754
$ws.Range("M122").Value = -7074.940000000001
$ws.Range("N122").Value = -46654

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2500
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 64747.46
$ws.Range("I136").Value = 87939.75
$ws.Range("J136").Value = 27639.8
$ws.Range("K136").Value = 263819.25
$ws.Range("L136").Value = 82919.39999999999
$ws.Range("M136").Value = -261269.25
$ws.Range("N136").Value = -88019.39999999999
